$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (row 33), pushing the
# existing rows 33:151 down to 35:153. Excel's native Insert semantics take
# care of shifting all the existing cell values/styles down by two rows,
# which is exactly the weekly "new week pushed in / oldest week pushed out"
# update this dataset receives.
$ws.Rows("33:34").Insert()

# New week's data (Primera / Segunda) for 2021-08-06.
$ws.Range("A33").Value = 1
$ws.Range("B33").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C33").Value = "Arica y Parinacota"
$ws.Range("D33").Value = 44414
$ws.Range("E33").Value = 15
$ws.Range("F33").Value = 100112043
$ws.Range("G33").Value = "Pepino ensalada"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 120
$ws.Range("K33").Value = 13000
$ws.Range("L33").Value = 14000
$ws.Range("M33").Value = 13500
$ws.Range("N33").Value = "$/caja 70 unidades"
$ws.Range("O33").Value = "Región de Arica y Parinacota"
$ws.Range("P33").Value = 193
$ws.Range("Q33").Value = 70
$ws.Range("R33").Value = "Hortaliza"

$ws.Range("A34").Value = 1
$ws.Range("B34").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C34").Value = "Arica y Parinacota"
$ws.Range("D34").Value = 44414
$ws.Range("E34").Value = 15
$ws.Range("F34").Value = 100112043
$ws.Range("G34").Value = "Pepino ensalada"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Segunda"
$ws.Range("J34").Value = 150
$ws.Range("K34").Value = 8000
$ws.Range("L34").Value = 9000
$ws.Range("M34").Value = 8500
$ws.Range("N34").Value = "$/caja 100 unidades"
$ws.Range("O34").Value = "Región de Arica y Parinacota"
$ws.Range("P34").Value = 85
$ws.Range("Q34").Value = 100
$ws.Range("R34").Value = "Hortaliza"
